# B6-PowerPoint.pptx edit replay
#
# 1) Re-style the three tables (slides 14, 15, 16) from the bespoke
#    "Table_0" style {1C0EADD3-EF24-4EAA-AF87-E6F8B3ED9C14} to the
#    built-in table style {C38530F9-00BD-4E73-9AD8-3D32A4C11A2F}.
# 2) Re-apply the deck's design from "Integral" (Red Violet) to the
#    default "Office Theme" colour scheme, reproducing the swap of the
#    rendered theme colours that the Design-gallery re-pick produced.

$p = $ppt.ActivePresentation

# --- 1. Table styles --------------------------------------------------
$newTableStyleId = "{C38530F9-00BD-4E73-9AD8-3D32A4C11A2F}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyleId)
    }
}

# --- 2. Theme colours ---------------------------------------------------
# Office Theme colour scheme values (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), applied in MsoThemeColorSchemeIndex order.
$officeThemeColors = @(
    0x000000,  # Text 1 / dk1
    0xFFFFFF,  # Background 1 / lt1
    0x44546A,  # Text 2 / dk2
    0xE7E6E6,  # Background 2 / lt2
    0x5B9BD5,  # Accent 1
    0xED7D31,  # Accent 2
    0xA5A5A5,  # Accent 3
    0xFFC000,  # Accent 4
    0x4472C4,  # Accent 5
    0x70AD47,  # Accent 6
    0x0563C1,  # Hyperlink
    0x954F72   # Followed Hyperlink
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $rgbHex = $officeThemeColors[$i - 1]
    # COM RGB() packs colour as 0x00BBGGRR, so byte-swap R and B.
    $r = [math]::Floor($rgbHex / 0x10000) % 0x100
    $g = [math]::Floor($rgbHex / 0x100) % 0x100
    $b = $rgbHex % 0x100
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $colorScheme.Item($i).RGB = $bgr
}
